$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns stay text (avoid Excel auto-numeric coercion
# that would strip trailing zeros / use scientific notation).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.478.90'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '1.911.16'
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '239.73'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = '0.9991'
$ws.Range('D7').Value = '0.4779'
$ws.Range('E7').Value = '  -2.78%  '
$ws.Range('D8').Value = '0.2848'
$ws.Range('E8').Value = '  -3.22%  '
$ws.Range('D9').Value = '0.06711'
$ws.Range('E9').Value = '  -3.03%  '
$ws.Range('D10').Value = '19.48'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('D11').Value = '103.41'
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').Value = '0.07765'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '1.917.02'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = '5.193'
$ws.Range('E14').Value = '  -3.42%  '
$ws.Range('D15').Value = '0.6708'
$ws.Range('E15').Value = '  -4.90%  '
$ws.Range('D16').Value = '276.00'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '30.506.86'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = '0.9990'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '0.000007493'
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').Value = '  -3.80%  '
$ws.Range('D21').Value = '5.382'
$ws.Range('E21').Value = '  -3.84%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '0.9989'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '6.300'
$ws.Range('E23').Value = '  -3.79%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '9.355'
$ws.Range('E24').Value = '  -4.79%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '167.19'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '19.25'
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.083'
$ws.Range('E27').Value = '  -3.52%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.388'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '0.09992'
$ws.Range('E29').Value = '  -4.22%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '4.593'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.511'
$ws.Range('E31').Value = '  -3.18%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.263'
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.04718'
$ws.Range('E33').Value = '  -3.59%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.7292'
$ws.Range('E34').Value = '  -3.88%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.117'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.713'
$ws.Range('E36').Value = '  -0.84%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01906'
$ws.Range('E37').Value = '  -5.01%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.608'
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '6.342'
$ws.Range('E39').Value = '  -2.05%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '74.19'
$ws.Range('E40').Value = '  -5.69%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '1.963'
$ws.Range('E41').Value = '  -6.38%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8619'
$ws.Range('E42').Value = '  -5.75%  '
$ws.Range('D43').Value = '106.56'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4265'
$ws.Range('E44').Value = '  -3.94%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9987'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.428'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '953.28'
$ws.Range('E47').Value = '  -4.17%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1211'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '34.63'
$ws.Range('E49').Value = '  -4.06%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05795'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '8.746'
$ws.Range('E51').Value = '  -4.58%  '
